# Update "想去人数" (want-to-go count) values for two events that are
# duplicated across the "展览" and "全部类型" worksheets.
#   F4: 31  -> 32
#   F5: 267 -> 268

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 32
    $ws.Range("F5").Value = 268
}
